# Generate Report for Handback
# Adds a new handback record (6419ced1-b836-49b0-bd9d-e5c48dc6a783) as row 4
# on the "Overview", "zh-cn" and "de-de" worksheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newId      = "6419ced1-b836-49b0-bd9d-e5c48dc6a783"
$newHash    = "a4882b8a1033b12df8e012010aaef65087f13a0c"
$statusSync = "Handed back: in sync with en-US"

$mdName     = "$newId.md"
$zhXlfName  = "$newId.$newHash.zh-cn.xlf"
$deXlfName  = "$newId.$newHash.de-de.xlf"

# ---------------------------------------------------------------------------
# Re-assert the existing content first (sheet by sheet, row-major, in the
# same order the data was originally authored), then append the new row
# for the handed-back file.
# ---------------------------------------------------------------------------

# ---- Sheet "Overview" (A:C) ----
$wsOverview.Range("A1").Value = "File Name"
$wsOverview.Range("B1").Value = "zh-cn"
$wsOverview.Range("C1").Value = "de-de"

$wsOverview.Range("A2").Value = "0e5c8e81-b798-489a-8fb1-1a0450417ae6.md"
$wsOverview.Range("B2").Value = "Handed back: not in sync with en-US"
$wsOverview.Range("C2").Value = "Handed back: not in sync with en-US"

$wsOverview.Range("A3").Value = "11918b85-a760-496a-a70e-715dd9214897.md"
$wsOverview.Range("B3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("C3").Value = "Handed back: in sync with en-US"

$wsOverview.Range("A4").Value = $mdName
$wsOverview.Range("B4").Value = $statusSync
$wsOverview.Range("C4").Value = $statusSync

$h = $wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/ee8267a350f3b767006f915bac959430401357dc/e2e/$mdName", [Type]::Missing, [Type]::Missing, $mdName)
$wsOverview.Range("A4").Style = "HyperLink"

# ---- Sheet "zh-cn" (A:I) ----
$wsZhCn.Range("A1").Value = "Source File Name"
$wsZhCn.Range("B1").Value = "Status"
$wsZhCn.Range("C1").Value = "Correspond Handoff File"
$wsZhCn.Range("D1").Value = "Correspond Handoff Datetime"
$wsZhCn.Range("E1").Value = "Target File"
$wsZhCn.Range("F1").Value = "Correspond Handback File"
$wsZhCn.Range("G1").Value = "Correspond Handback DateTime"
$wsZhCn.Range("H1").Value = "Handoff Reason"
$wsZhCn.Range("I1").Value = "Dependency From"

$wsZhCn.Range("A2").Value = "0e5c8e81-b798-489a-8fb1-1a0450417ae6.md"
$wsZhCn.Range("B2").Value = "Handed back: not in sync with en-US"
$wsZhCn.Range("C2").Value = "0e5c8e81-b798-489a-8fb1-1a0450417ae6.11775246a1e0f778a57d015760652fcb68fa9064.zh-cn.xlf"
$wsZhCn.Range("D2").Value = "2016-02-19 06:56:34"
$wsZhCn.Range("E2").Value = "0e5c8e81-b798-489a-8fb1-1a0450417ae6.md"
$wsZhCn.Range("F2").Value = "0e5c8e81-b798-489a-8fb1-1a0450417ae6.11775246a1e0f778a57d015760652fcb68fa9064.zh-cn.xlf"
$wsZhCn.Range("G2").Value = "2016-02-19 06:57:54"
$wsZhCn.Range("H2").Value = "Include"

$wsZhCn.Range("A3").Value = "11918b85-a760-496a-a70e-715dd9214897.md"
$wsZhCn.Range("B3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C3").Value = "11918b85-a760-496a-a70e-715dd9214897.590f408322e42aa4a556f8e8b2c5586ef687bea2.zh-cn.xlf"
$wsZhCn.Range("D3").Value = "2016-02-19 07:02:39"
$wsZhCn.Range("E3").Value = "11918b85-a760-496a-a70e-715dd9214897.md"
$wsZhCn.Range("F3").Value = "11918b85-a760-496a-a70e-715dd9214897.590f408322e42aa4a556f8e8b2c5586ef687bea2.zh-cn.xlf"
$wsZhCn.Range("G3").Value = "2016-02-19 07:03:18"
$wsZhCn.Range("H3").Value = "Include"

$wsZhCn.Range("A4").Value = $mdName
$wsZhCn.Range("B4").Value = $statusSync
$wsZhCn.Range("C4").Value = $zhXlfName
$wsZhCn.Range("D4").Value = "2016-02-19 07:05:30"
$wsZhCn.Range("E4").Value = $mdName
$wsZhCn.Range("F4").Value = $zhXlfName
$wsZhCn.Range("G4").Value = "2016-02-19 07:06:13"
$wsZhCn.Range("H4").Value = "Include"

$h = $wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/ee8267a350f3b767006f915bac959430401357dc/e2e/$mdName", [Type]::Missing, [Type]::Missing, $mdName)
$h = $wsZhCn.Hyperlinks.Add($wsZhCn.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/50d9b72ef4ed396d21db423dd7f4f61fb3c6dc28/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/$zhXlfName", [Type]::Missing, [Type]::Missing, $zhXlfName)
$h = $wsZhCn.Hyperlinks.Add($wsZhCn.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/4a081b3654f30f87cc81bee416c6c74f14bf4fd3/e2e/$mdName", [Type]::Missing, [Type]::Missing, $mdName)
$h = $wsZhCn.Hyperlinks.Add($wsZhCn.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/ec40c6b5a6041fd3a5c76f3cc454909db71ffbf5/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/$zhXlfName", [Type]::Missing, [Type]::Missing, $zhXlfName)

$wsZhCn.Range("A4").Style = "HyperLink"
$wsZhCn.Range("C4").Style = "HyperLink"
$wsZhCn.Range("E4").Style = "HyperLink"
$wsZhCn.Range("F4").Style = "HyperLink"
$wsZhCn.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---- Sheet "de-de" (A:I) ----
$wsDeDe.Range("A1").Value = "Source File Name"
$wsDeDe.Range("B1").Value = "Status"
$wsDeDe.Range("C1").Value = "Correspond Handoff File"
$wsDeDe.Range("D1").Value = "Correspond Handoff Datetime"
$wsDeDe.Range("E1").Value = "Target File"
$wsDeDe.Range("F1").Value = "Correspond Handback File"
$wsDeDe.Range("G1").Value = "Correspond Handback DateTime"
$wsDeDe.Range("H1").Value = "Handoff Reason"
$wsDeDe.Range("I1").Value = "Dependency From"

$wsDeDe.Range("A2").Value = "0e5c8e81-b798-489a-8fb1-1a0450417ae6.md"
$wsDeDe.Range("B2").Value = "Handed back: not in sync with en-US"
$wsDeDe.Range("C2").Value = "0e5c8e81-b798-489a-8fb1-1a0450417ae6.11775246a1e0f778a57d015760652fcb68fa9064.de-de.xlf"
$wsDeDe.Range("D2").Value = "2016-02-19 06:56:44"
$wsDeDe.Range("E2").Value = "0e5c8e81-b798-489a-8fb1-1a0450417ae6.md"
$wsDeDe.Range("F2").Value = "0e5c8e81-b798-489a-8fb1-1a0450417ae6.11775246a1e0f778a57d015760652fcb68fa9064.de-de.xlf"
$wsDeDe.Range("G2").Value = "2016-02-19 06:58:10"
$wsDeDe.Range("H2").Value = "Include"

$wsDeDe.Range("A3").Value = "11918b85-a760-496a-a70e-715dd9214897.md"
$wsDeDe.Range("B3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C3").Value = "11918b85-a760-496a-a70e-715dd9214897.590f408322e42aa4a556f8e8b2c5586ef687bea2.de-de.xlf"
$wsDeDe.Range("D3").Value = "2016-02-19 07:02:48"
$wsDeDe.Range("E3").Value = "11918b85-a760-496a-a70e-715dd9214897.md"
$wsDeDe.Range("F3").Value = "11918b85-a760-496a-a70e-715dd9214897.590f408322e42aa4a556f8e8b2c5586ef687bea2.de-de.xlf"
$wsDeDe.Range("G3").Value = "2016-02-19 07:03:35"
$wsDeDe.Range("H3").Value = "Include"

$wsDeDe.Range("A4").Value = $mdName
$wsDeDe.Range("B4").Value = $statusSync
$wsDeDe.Range("C4").Value = $deXlfName
$wsDeDe.Range("D4").Value = "2016-02-19 07:05:40"
$wsDeDe.Range("E4").Value = $mdName
$wsDeDe.Range("F4").Value = $deXlfName
$wsDeDe.Range("G4").Value = "2016-02-19 07:06:30"
$wsDeDe.Range("H4").Value = "Include"

$h = $wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/ee8267a350f3b767006f915bac959430401357dc/e2e/$mdName", [Type]::Missing, [Type]::Missing, $mdName)
$h = $wsDeDe.Hyperlinks.Add($wsDeDe.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c94d76291030168f904cf6f3870ac069596f0469/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/$deXlfName", [Type]::Missing, [Type]::Missing, $deXlfName)
$h = $wsDeDe.Hyperlinks.Add($wsDeDe.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/b1ac34266f8a2354c74c6de7f576d3dd3d741d98/e2e/$mdName", [Type]::Missing, [Type]::Missing, $mdName)
$h = $wsDeDe.Hyperlinks.Add($wsDeDe.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/4b9b68b0d2af03f5ef38b127f3ee40be183ea45a/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/$deXlfName", [Type]::Missing, [Type]::Missing, $deXlfName)

$wsDeDe.Range("A4").Style = "HyperLink"
$wsDeDe.Range("C4").Style = "HyperLink"
$wsDeDe.Range("E4").Style = "HyperLink"
$wsDeDe.Range("F4").Style = "HyperLink"
$wsDeDe.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

Write-Output "Report generated for handback $newId"
